$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 6).Value = 65   # F2
$ws.Cells.Item(2, 7).Value = 90   # G2
$ws.Cells.Item(2, 8).Value = 1.07   # H2
$ws.Cells.Item(2, 9).Value = 1.09   # I2
$ws.Cells.Item(2, 10).Value = 15   # J2
$ws.Cells.Item(2, 11).Value = 18.5   # K2
$ws.Cells.Item(2, 12).Value = 0   # L2
$ws.Cells.Item(2, 13).Value = 0   # M2
$ws.Cells.Item(2, 14).Value = 0   # N2
$ws.Cells.Item(2, 15).Value = 0   # O2
$ws.Cells.Item(2, 16).Value = 0   # P2
$ws.Cells.Item(2, 17).Value = 0   # Q2
$ws.Cells.Item(2, 18).Value = 11   # R2
$ws.Cells.Item(2, 19).Value = 1.09   # S2
$ws.Cells.Item(2, 20).Value = 0   # T2
$ws.Cells.Item(2, 21).Value = 0   # U2
$ws.Cells.Item(2, 22).Value = 12   # V2
$ws.Cells.Item(2, 23).Value = 1.01   # W2
$ws.Cells.Item(2, 24).Value = 1000   # X2
$ws.Cells.Item(2, 25).Value = 1000   # Y2
$ws.Cells.Item(2, 26).Value = 1000   # Z2
$ws.Cells.Item(2, 27).Value = 1000   # AA2
$ws.Cells.Item(2, 28).Value = 1000   # AB2
$ws.Cells.Item(2, 29).Value = 1000   # AC2
$ws.Cells.Item(2, 30).Value = 11.5   # AD2
$ws.Cells.Item(2, 31).Value = 8   # AE2
$ws.Cells.Item(2, 33).Value = 1000   # AG2
$ws.Cells.Item(2, 34).Value = 540   # AH2
$ws.Cells.Item(2, 35).Value = 990   # AI2
$ws.Cells.Item(2, 41).Value = 3.2   # AO2
# Row 3
$ws.Cells.Item(3, 6).Value = 2.24   # F3
$ws.Cells.Item(3, 7).Value = 2.32   # G3
$ws.Cells.Item(3, 8).Value = 4.2   # H3
$ws.Cells.Item(3, 9).Value = 4.7   # I3
$ws.Cells.Item(3, 10).Value = 2.96   # J3
$ws.Cells.Item(3, 11).Value = 3   # K3
$ws.Cells.Item(3, 14).Value = 2.36   # N3
$ws.Cells.Item(3, 15).Value = 1.7   # O3
$ws.Cells.Item(3, 16).Value = 1.42   # P3
$ws.Cells.Item(3, 17).Value = 3.2   # Q3
$ws.Cells.Item(3, 18).Value = 1.14   # R3
$ws.Cells.Item(3, 20).Value = 2.5   # T3
$ws.Cells.Item(3, 21).Value = 1.57   # U3
$ws.Cells.Item(3, 22).Value = 1.28   # V3
$ws.Cells.Item(3, 23).Value = 1.76   # W3
$ws.Cells.Item(3, 24).Value = 7.4   # X3
$ws.Cells.Item(3, 25).Value = 10.5   # Y3
$ws.Cells.Item(3, 26).Value = 42   # Z3
$ws.Cells.Item(3, 27).Value = 1000   # AA3
$ws.Cells.Item(3, 28).Value = 5.9   # AB3
$ws.Cells.Item(3, 30).Value = 22   # AD3
$ws.Cells.Item(3, 31).Value = 240   # AE3
$ws.Cells.Item(3, 32).Value = 12.5   # AF3
$ws.Cells.Item(3, 33).Value = 13.5   # AG3
$ws.Cells.Item(3, 34).Value = 40   # AH3
$ws.Cells.Item(3, 35).Value = 1000   # AI3
$ws.Cells.Item(3, 36).Value = 34   # AJ3
$ws.Cells.Item(3, 37).Value = 65   # AK3
$ws.Cells.Item(3, 40).Value = 120   # AN3
# Row 4
$ws.Cells.Item(4, 6).Value = 2.02   # F4
$ws.Cells.Item(4, 8).Value = 3.45   # H4
$ws.Cells.Item(4, 9).Value = 3.7   # I4
$ws.Cells.Item(4, 16).Value = 2.82   # P4
$ws.Cells.Item(4, 17).Value = 1.52   # Q4
$ws.Cells.Item(4, 18).Value = 1.74   # R4
$ws.Cells.Item(4, 19).Value = 2.3   # S4
$ws.Cells.Item(4, 20).Value = 1.52   # T4
$ws.Cells.Item(4, 21).Value = 2.74   # U4
$ws.Cells.Item(4, 22).Value = 1.37   # V4
$ws.Cells.Item(4, 25).Value = 22   # Y4
$ws.Cells.Item(4, 28).Value = 16   # AB4
$ws.Cells.Item(4, 29).Value = 10.5   # AC4
$ws.Cells.Item(4, 30).Value = 16   # AD4
$ws.Cells.Item(4, 31).Value = 36   # AE4
$ws.Cells.Item(4, 32).Value = 16.5   # AF4
$ws.Cells.Item(4, 34).Value = 14.5   # AH4
$ws.Cells.Item(4, 36).Value = 26   # AJ4
$ws.Cells.Item(4, 37).Value = 18   # AK4
$ws.Cells.Item(4, 38).Value = 26   # AL4
$ws.Cells.Item(4, 39).Value = 55   # AM4
$ws.Cells.Item(4, 40).Value = 8.800000000000001   # AN4
$ws.Cells.Item(4, 41).Value = 22   # AO4
# Row 5
$ws.Cells.Item(5, 6).Value = 2.88   # F5
$ws.Cells.Item(5, 7).Value = 3.2   # G5
$ws.Cells.Item(5, 8).Value = 3.05   # H5
$ws.Cells.Item(5, 10).Value = 2.6   # J5
$ws.Cells.Item(5, 13).Value = 1.18   # M5
$ws.Cells.Item(5, 14).Value = 2.3   # N5
$ws.Cells.Item(5, 15).Value = 1.7   # O5
$ws.Cells.Item(5, 17).Value = 3.25   # Q5
$ws.Cells.Item(5, 19).Value = 7   # S5
$ws.Cells.Item(5, 21).Value = 1.57   # U5
$ws.Cells.Item(5, 22).Value = 1.41   # V5
# Row 6
$ws.Cells.Item(6, 6).Value = 1.38   # F6
$ws.Cells.Item(6, 7).Value = 1.4   # G6
$ws.Cells.Item(6, 8).Value = 9.4   # H6
$ws.Cells.Item(6, 9).Value = 10.5   # I6
$ws.Cells.Item(6, 10).Value = 5.5   # J6
$ws.Cells.Item(6, 11).Value = 5.9   # K6
$ws.Cells.Item(6, 12).Value = 1.34   # L6
$ws.Cells.Item(6, 14).Value = 4.8   # N6
$ws.Cells.Item(6, 16).Value = 2.32   # P6
$ws.Cells.Item(6, 17).Value = 1.71   # Q6
$ws.Cells.Item(6, 18).Value = 1.5   # R6
$ws.Cells.Item(6, 19).Value = 2.84   # S6
$ws.Cells.Item(6, 20).Value = 1.98   # T6
$ws.Cells.Item(6, 21).Value = 1.86   # U6
$ws.Cells.Item(6, 22).Value = 1.1   # V6
$ws.Cells.Item(6, 23).Value = 3.45   # W6
$ws.Cells.Item(6, 24).Value = 24   # X6
$ws.Cells.Item(6, 25).Value = 36   # Y6
$ws.Cells.Item(6, 26).Value = 95   # Z6
$ws.Cells.Item(6, 27).Value = 340   # AA6
$ws.Cells.Item(6, 29).Value = 12.5   # AC6
$ws.Cells.Item(6, 30).Value = 38   # AD6
$ws.Cells.Item(6, 31).Value = 160   # AE6
$ws.Cells.Item(6, 32).Value = 8.4   # AF6
$ws.Cells.Item(6, 33).Value = 10.5   # AG6
$ws.Cells.Item(6, 34).Value = 29   # AH6
$ws.Cells.Item(6, 35).Value = 290   # AI6
$ws.Cells.Item(6, 36).Value = 11.5   # AJ6
$ws.Cells.Item(6, 38).Value = 36   # AL6
$ws.Cells.Item(6, 39).Value = 160   # AM6
$ws.Cells.Item(6, 40).Value = 6.2   # AN6
$ws.Cells.Item(6, 41).Value = 190   # AO6
# Row 7
$ws.Cells.Item(7, 6).Value = 2.74   # F7
$ws.Cells.Item(7, 7).Value = 2.88   # G7
$ws.Cells.Item(7, 8).Value = 2.7   # H7
$ws.Cells.Item(7, 9).Value = 2.84   # I7
$ws.Cells.Item(7, 11).Value = 3.65   # K7
$ws.Cells.Item(7, 12).Value = 1.44   # L7
$ws.Cells.Item(7, 14).Value = 3.65   # N7
$ws.Cells.Item(7, 16).Value = 1.88   # P7
$ws.Cells.Item(7, 18).Value = 1.33   # R7
$ws.Cells.Item(7, 19).Value = 3.75   # S7
$ws.Cells.Item(7, 20).Value = 1.79   # T7
$ws.Cells.Item(7, 22).Value = 1.54   # V7
$ws.Cells.Item(7, 23).Value = 1.53   # W7
$ws.Cells.Item(7, 24).Value = 13.5   # X7
$ws.Cells.Item(7, 25).Value = 11   # Y7
$ws.Cells.Item(7, 26).Value = 46   # Z7
$ws.Cells.Item(7, 27).Value = 900   # AA7
$ws.Cells.Item(7, 30).Value = 12.5   # AD7
$ws.Cells.Item(7, 32).Value = 29   # AF7
$ws.Cells.Item(7, 33).Value = 13   # AG7
$ws.Cells.Item(7, 35).Value = 980   # AI7
$ws.Cells.Item(7, 36).Value = 130   # AJ7
$ws.Cells.Item(7, 37).Value = 44   # AK7
$ws.Cells.Item(7, 38).Value = 290   # AL7
$ws.Cells.Item(7, 40).Value = 46   # AN7
$ws.Cells.Item(7, 41).Value = 70   # AO7
# Row 8
$ws.Cells.Item(8, 7).Value = 1.5   # G8
$ws.Cells.Item(8, 9).Value = 13   # I8
$ws.Cells.Item(8, 11).Value = 4.6   # K8
$ws.Cells.Item(8, 23).Value = 2.98   # W8
$ws.Cells.Item(8, 25).Value = 42   # Y8
$ws.Cells.Item(8, 26).Value = 120   # Z8
$ws.Cells.Item(8, 28).Value = 5.3   # AB8
$ws.Cells.Item(8, 29).Value = 12   # AC8
$ws.Cells.Item(8, 35).Value = 1000   # AI8
$ws.Cells.Item(8, 36).Value = 12   # AJ8
$ws.Cells.Item(8, 37).Value = 80   # AK8
$ws.Cells.Item(8, 38).Value = 110   # AL8
$ws.Cells.Item(8, 40).Value = 1000   # AN8
# Row 9
$ws.Cells.Item(9, 9).Value = 7.4   # I9
$ws.Cells.Item(9, 15).Value = 1.47   # O9
$ws.Cells.Item(9, 18).Value = 1.26   # R9
$ws.Cells.Item(9, 20).Value = 2.36   # T9
$ws.Cells.Item(9, 21).Value = 1.71   # U9
$ws.Cells.Item(9, 25).Value = 18   # Y9
$ws.Cells.Item(9, 28).Value = 6.2   # AB9
$ws.Cells.Item(9, 33).Value = 11   # AG9
$ws.Cells.Item(9, 35).Value = 140   # AI9
$ws.Cells.Item(9, 41).Value = 230   # AO9
# Row 10
$ws.Cells.Item(10, 9).Value = 38   # I10
$ws.Cells.Item(10, 15).Value = 1.12   # O10
$ws.Cells.Item(10, 16).Value = 3.25   # P10
$ws.Cells.Item(10, 19).Value = 2.06   # S10
$ws.Cells.Item(10, 20).Value = 2.86   # T10
$ws.Cells.Item(10, 23).Value = 9.800000000000001   # W10
$ws.Cells.Item(10, 25).Value = 110   # Y10
$ws.Cells.Item(10, 26).Value = 480   # Z10
$ws.Cells.Item(10, 29).Value = 990   # AC10
$ws.Cells.Item(10, 30).Value = 1000   # AD10
$ws.Cells.Item(10, 32).Value = 8   # AF10
$ws.Cells.Item(10, 34).Value = 1000   # AH10
$ws.Cells.Item(10, 36).Value = 7.2   # AJ10
$ws.Cells.Item(10, 37).Value = 1000   # AK10
$ws.Cells.Item(10, 38).Value = 1000   # AL10
$ws.Cells.Item(10, 40).Value = 2.9   # AN10
# Row 11
$ws.Cells.Item(11, 6).Value = 5.6   # F11
$ws.Cells.Item(11, 7).Value = 9.4   # G11
$ws.Cells.Item(11, 8).Value = 1.42   # H11
$ws.Cells.Item(11, 9).Value = 1.56   # I11
$ws.Cells.Item(11, 10).Value = 4.1   # J11
$ws.Cells.Item(11, 11).Value = 5.8   # K11
$ws.Cells.Item(11, 12).Value = 1.25   # L11
$ws.Cells.Item(11, 13).Value = 1.03   # M11
$ws.Cells.Item(11, 14).Value = 4.1   # N11
$ws.Cells.Item(11, 15).Value = 1.18   # O11
$ws.Cells.Item(11, 16).Value = 2.22   # P11
$ws.Cells.Item(11, 17).Value = 1.54   # Q11
$ws.Cells.Item(11, 18).Value = 1.53   # R11
$ws.Cells.Item(11, 19).Value = 2.22   # S11
$ws.Cells.Item(11, 20).Value = 1.71   # T11
$ws.Cells.Item(11, 21).Value = 2   # U11
$ws.Cells.Item(11, 22).Value = 2.78   # V11
$ws.Cells.Item(11, 23).Value = 1.13   # W11
$ws.Cells.Item(11, 25).Value = 13   # Y11
$ws.Cells.Item(11, 26).Value = 12.5   # Z11
$ws.Cells.Item(11, 27).Value = 16.5   # AA11
$ws.Cells.Item(11, 29).Value = 14.5   # AC11
$ws.Cells.Item(11, 30).Value = 12.5   # AD11
$ws.Cells.Item(11, 31).Value = 18   # AE11
$ws.Cells.Item(11, 41).Value = 7   # AO11

Write-Host "Applied 215 cell updates."
